$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the template row (row 48) down through the new rows
$ws.Range("A48:D48").Copy($ws.Range("A49:D62"))

# Fill in the new time-log entries (rows 49-62)
$ws.Range("A49").Value = 44676.625
$ws.Range("B49").Value = 44676.6875
$ws.Range("C49").Formula = "=B49-A49"
$ws.Range("D49").Formula = "=C49+D48"
$ws.Range("E49").Value = "TUI Implementation"

$ws.Range("A50").Value = 44677.75
$ws.Range("B50").Value = 44677.8125
$ws.Range("C50").Formula = "=B50-A50"
$ws.Range("D50").Formula = "=C50+D49"
$ws.Range("E50").Value = "TUI Implementation"

$ws.Range("A51").Value = 44678.583333333336
$ws.Range("B51").Value = 44678.604166666664
$ws.Range("C51").Formula = "=B51-A51"
$ws.Range("D51").Formula = "=C51+D50"
$ws.Range("E51").Value = "Meeting with Phil to debug window size error"

$ws.Range("A52").Value = 44676.604166666664
$ws.Range("B52").Value = 44676.625
$ws.Range("C52").Formula = "=B52-A52"
$ws.Range("D52").Formula = "=C52+D51"
$ws.Range("E52").Value = "Client Meeting"

$ws.Range("A53").Value = 44680.6875
$ws.Range("B53").Value = 44680.75
$ws.Range("C53").Formula = "=B53-A53"
$ws.Range("D53").Formula = "=C53+D52"
$ws.Range("E53").Value = "Placeholder Network display and taking user iput for network selection"

$ws.Range("A54").Value = 44683.791666666664
$ws.Range("B54").Value = 44683.822916666664
$ws.Range("C54").Formula = "=B54-A54"
$ws.Range("D54").Formula = "=C54+D53"
$ws.Range("E54").Value = "User input for network selection and brainstorming names"

$ws.Range("A55").Value = 44683.604166666664
$ws.Range("B55").Value = 44683.625
$ws.Range("C55").Formula = "=B55-A55"
$ws.Range("D55").Formula = "=C55+D54"
$ws.Range("E55").Value = "Client meeting"

$ws.Range("A56").Value = 44686.697916666664
$ws.Range("B56").Value = 44686.760416666664
$ws.Range("C56").Formula = "=B56-A56"
$ws.Range("D56").Formula = "=C56+D55"
$ws.Range("E56").Value = "Working on TUI subwindow and transfering work to master branch"

$ws.Range("A57").Value = 44690.770833333336
$ws.Range("B57").Value = 44690.8125
$ws.Range("C57").Formula = "=B57-A57"
$ws.Range("D57").Formula = "=C57+D56"
$ws.Range("E57").Value = "Work on BDS Makefile"

$ws.Range("A58").Value = 44692.604166666664
$ws.Range("B58").Value = 44692.625
$ws.Range("C58").Formula = "=B58-A58"
$ws.Range("D58").Formula = "=C58+D57"
$ws.Range("E58").Value = "Client Meeting"

$ws.Range("A59").Value = 44694.729166666664
$ws.Range("B59").Value = 44694.802083333336
$ws.Range("C59").Formula = "=B59-A59"
$ws.Range("D59").Formula = "=C59+D58"
$ws.Range("E59").Value = "Work on BSD Makefile and Subwindow configuration"

$ws.Range("A60").Value = 44697.041666666664
$ws.Range("B60").Value = 44697.104166666664
$ws.Range("C60").Formula = "=B60-A60"
$ws.Range("D60").Formula = "=C60+D59"
$ws.Range("E60").Value = "Implementation of Network Display"

$ws.Range("A61").Value = 44697.604166666664
$ws.Range("B61").Value = 44697.625
$ws.Range("C61").Formula = "=B61-A61"
$ws.Range("D61").Formula = "=C61+D60"
$ws.Range("E61").Value = "Client Meeting"

$ws.Range("A62").Value = 44699.458333333336
$ws.Range("B62").Value = 44699.5
$ws.Range("C62").Formula = "=B62-A62"
$ws.Range("D62").Formula = "=C62+D61"
$ws.Range("E62").Value = "Milestone 9 and time log updates"

# Update the active selection to reflect the next empty entry row, matching
# where the cursor was left after adding the new log entries.
$ws.Range("E63").Select()
